$wb = $excel.ActiveWorkbook

# Target tab order: 总计, 2022-Q2, 2022-Q1 (was: 2022-Q1, 2022-Q2, 总计)
# i.e. "总计" moves to the front and "2022-Q1" moves to the back;
# "2022-Q2" stays put in the middle. Each sheet keeps its own
# content/formatting - only the tab order (and which tab is active) changes.

# 1) Move "总计" to be the first sheet.
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal.Move($wsQ1)

# 2) Move "2022-Q1" to be right after "2022-Q2" (i.e. to the end).
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ1.Move($null, $wsQ2)

# 3) "2022-Q1" was (and remains) the selected/active tab.
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Activate()
